$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of G2 and G3 (pompe Back Left <-> pompe Back Right)
$g2 = $ws.Range("G2").Value
$g3 = $ws.Range("G3").Value
$ws.Range("G2").Value = $g3
$ws.Range("G3").Value = $g2

# Update the view: scroll so column B is the leftmost visible column,
# and move the active selection to G3
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("G3").Select()
